# Update cryptos list values per diff (price/volume refresh, row13<->14 and row45<->46 swaps)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.450.73"
$ws.Range("E2").Value = "  -1.15%  "
$ws.Range("D3").Value = "'1.830.29"
$ws.Range("E3").Value = "  -2.22%  "
$ws.Range("E4").Value = "  -0.84%  "
$ws.Range("D5").Value = "'331.57"
$ws.Range("E5").Value = "  -0.95%  "
$ws.Range("E6").Value = "  -0.70%  "
$ws.Range("D7").Value = "'0.4579"
$ws.Range("E8").Value = "  -3.15%  "
$ws.Range("D9").Value = "'46.49"
$ws.Range("E9").Value = "  +0.98%  "
$ws.Range("E10").Value = "  -1.35%  "
$ws.Range("D11").Value = "'0.9679"
$ws.Range("E11").Value = "  -3.99%  "
$ws.Range("E12").Value = "  -3.96%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'5.880"
$ws.Range("E13").Value = "  -2.41%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "'1.824.21"
$ws.Range("E14").Value = "  -2.60%  "
$ws.Range("D15").Value = "'7.069"
$ws.Range("E15").Value = "  -2.94%  "
$ws.Range("D16").Value = "'1.004"
$ws.Range("D17").Value = "'89.40"
$ws.Range("E17").Value = "  +0.59%  "
$ws.Range("D18").Value = "'0.06606"
$ws.Range("E18").Value = "  -2.06%  "
$ws.Range("E19").Value = "  -1.85%  "
$ws.Range("D20").Value = "'17.12"
$ws.Range("E20").Value = "  -0.66%  "
$ws.Range("E21").Value = "  -0.73%  "
$ws.Range("D22").Value = "'27.440.13"
$ws.Range("E22").Value = "  -1.02%  "
$ws.Range("E23").Value = "  -2.72%  "
$ws.Range("D24").Value = "'10.80"
$ws.Range("E24").Value = "  -1.44%  "
$ws.Range("D25").Value = "'2.292"
$ws.Range("E25").Value = "  -1.05%  "
$ws.Range("D26").Value = "'2.055.10"
$ws.Range("E26").Value = "  -1.80%  "
$ws.Range("D27").Value = "'155.70"
$ws.Range("E27").Value = "  -2.08%  "
$ws.Range("D28").Value = "'19.38"
$ws.Range("E28").Value = "  -2.49%  "
$ws.Range("D29").Value = "'2.066"
$ws.Range("E29").Value = "  -4.99%  "
$ws.Range("D30").Value = "'5.295"
$ws.Range("E30").Value = "  -3.13%  "
$ws.Range("D31").Value = "'118.35"
$ws.Range("E31").Value = "  -3.04%  "
$ws.Range("D32").Value = "'0.9409"
$ws.Range("E32").Value = "  -4.29%  "
$ws.Range("D33").Value = "'0.09304"
$ws.Range("E33").Value = "  -2.08%  "
$ws.Range("E34").Value = "  -1.03%  "
$ws.Range("D35").Value = "'5.240"
$ws.Range("E35").Value = "  -1.70%  "
$ws.Range("D36").Value = "'1.328"
$ws.Range("E36").Value = "  -1.52%  "
$ws.Range("D37").Value = "'0.05921"
$ws.Range("E37").Value = "  -2.53%  "
$ws.Range("D38").Value = "'0.02177"
$ws.Range("E38").Value = "  -2.79%  "
$ws.Range("D39").Value = "'8.096"
$ws.Range("E39").Value = "  -3.27%  "
$ws.Range("D40").Value = "'1.150"
$ws.Range("E40").Value = "  -4.17%  "
$ws.Range("D41").Value = "'0.5776"
$ws.Range("E41").Value = "  -3.70%  "
$ws.Range("D42").Value = "'0.1828"
$ws.Range("E42").Value = "  -3.28%  "
$ws.Range("D43").Value = "'9.988"
$ws.Range("E43").Value = "  -3.52%  "
$ws.Range("D44").Value = "'1.268"
$ws.Range("E44").Value = "  +2.07%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'11.97"
$ws.Range("E45").Value = "  -1.90%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "'0.5449"
$ws.Range("E46").Value = "  -3.95%  "
$ws.Range("D47").Value = "'1.869"
$ws.Range("E47").Value = "  -3.27%  "
$ws.Range("D48").Value = "'110.76"
$ws.Range("E48").Value = "  -1.62%  "
$ws.Range("D49").Value = "'0.06596"
$ws.Range("E49").Value = "  -2.50%  "
$ws.Range("E50").Value = "  -0.84%  "
$ws.Range("E51").Value = "  -1.52%  "
